$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows with results that came in ---
$ws.Range("G34").Value = "Fallo"
$ws.Range("H34").Value = -1

$ws.Range("G80").Value = "Fallo"
$ws.Range("H80").Value = -1

$ws.Range("G82").Value = "Fallo"
$ws.Range("H82").Value = -1

$ws.Range("G88").Value = "Acierto"
$ws.Range("H88").Value = 2

# --- Append new pronostico rows (89-94) ---
# Column B holds dates formatted as plain text (e.g. "2025-09-03"), so force
# the text number format first to stop Excel from auto-converting them to
# date serial numbers.
$ws.Range("B89:B94").NumberFormat = "@"

$newRows = @(
    @(89, 14551800, "2025-09-03", "Dusan Lajovic", "Frederico Ferreira Silva", "Gana Frederico Ferreira Silva", 3.4),
    @(90, 14551798, "2025-09-03", "Ignacio Buse", "Pablo Llamas Ruiz", "Gana Ignacio Buse", 2.1),
    @(91, 14598777, "2025-09-03", "Maxime Janvier", "Mark Lajal", "Gana Maxime Janvier", 4.5),
    @(92, 14598719, "2025-09-03", "Kenny De Schepper", "Eliakim Coulibaly", "Gana Kenny De Schepper", 3.25),
    @(93, 14552652, "2025-09-03", "Abdullah Shelbayh", "Viktor Durasovic", "Gana Viktor Durasovic", 2.1),
    @(94, 14552665, "2025-09-03", "Nicolas Mejia", "Marek Gengel", "Gana Marek Gengel", 2.63)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
